$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy header style (bold, border, centered) from F1 into G1:H1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Updated precision values
$ws.Range("B2").Value = 0.3955935532374564
$ws.Range("D2").Value = 0.4321333824756292

# New data cells
$ws.Range("G2").Value = 0.1228190763666741
$ws.Range("H2").Value = 0.991
